$d = $word.ActiveDocument

$replacements = @(
    @("2025-05-19 Monday", "2025-05-20 Tuesday"),
    @("256×6=", "464×6="),
    @("937×4=", "619×9="),
    @("666×8=", "913×5="),
    @("681×5=", "555×3="),
    @("532×9=", "589×3="),
    @("852×2=", "243×4="),
    @("521×2=", "280×7="),
    @("221×9=", "281×6="),
    @("469×6=", "288×4="),
    @("897×4=", "853×9="),
    @("276×4=", "391×5="),
    @("219×6=", "886×7="),
    @("498×6=", "194×7="),
    @("190×7=", "817×5="),
    @("888×8=", "353×7="),
    @("669×4=", "186×6="),
    @("641×9=", "558×9="),
    @("251×9=", "706×6="),
    @("709×5=", "179×2="),
    @("650×6=", "485×6="),
    @("644×6=", "331×8="),
    @("218×9=", "611×3="),
    @("432×5=", "168×7="),
    @("236×2=", "386×7="),
    @("992×8=", "568×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
